$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap "Recorded By" column (G) order: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$recordedByRows = @(2,3,20,21,22,39,40,41,58,59,60,77,78,95,96,113,114,131,132,149,150,167,168,169,186,187,188,205,206,207)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# --- 2) Update summary metrics L7 (Missing Sessions) and L8 (Pending Sessions)
$ws.Range("L7").Value = 12
$ws.Range("L8").Value = 150

# --- 3) Update per-group statistics P/Q columns (rows 16-18 and 24-26)
$statRows = @(16,17,18,24,25,26)
foreach ($r in $statRows) {
    $ws.Range("P$r").Value = 1
    $ws.Range("Q$r").Value = 13
}

# --- 4) Rows whose last session flipped from "Pending" (yellow) to "Not Recorded" (pink)
# Reuse the exact formatting already used by existing "Not Recorded" rows (e.g. row 6)
# by copying its format, so the same pink fill style is applied consistently.
$notRecordedRows = @(24,43,62,171,190,209)
$formatSourceRow = 6
foreach ($r in $notRecordedRows) {
    $ws.Range("A$formatSourceRow`:I$formatSourceRow").Copy()
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
$excel.CutCopyMode = $false
